$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for new columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the style from an existing header cell (e.g. AC1) to the new header cells
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122) # xlPasteFormats

# Fill data rows 2-48 with team record values
for ($r = 2; $r -le 48; $r++) {
    $ws.Cells.Item($r, 30).Value = 74   # AD
    $ws.Cells.Item($r, 31).Value = 88   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
